$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 6, column C value (0.48958333333333331 -> 1)
$ws.Range("C6").Value = 1

# Add new row 7 data, copying formatting from row 6 so existing
# cell styles (date/time number formats) are reused instead of
# creating brand-new style entries.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 42376

$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").Value = 0

$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = 0.5541666666666667

$ws.Range("E7").Value = "Finished funcationality."

# Update active selection to E7
$ws.Range("E7").Select()
